$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure updated Price (D) cells remain plain text, matching the source data formatting
$priceCells = @("D2","D3","D4","D5","D6","D7","D11","D13","D15","D16","D18","D22","D27","D30","D31","D36","D40","D42","D43","D44","D47","D48","D50","D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated price (D) and volume/change (E) values row by row
$ws.Range("D2").Value = '41.783.36'
$ws.Range("E2").Value = '  +0.62%  '
$ws.Range("D3").Value = '2.478.73'
$ws.Range("E3").Value = '  +0.51%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '318.97'
$ws.Range("E5").Value = '  +1.58%  '
$ws.Range("D6").Value = '93.54'
$ws.Range("E6").Value = '  +2.42%  '
$ws.Range("D7").Value = '0.555'
$ws.Range("E7").Value = '  +0.98%  '
$ws.Range("E9").Value = '  +1.33%  '
$ws.Range("E10").Value = '  +11.29%  '
$ws.Range("D11").Value = '33.33'
$ws.Range("E11").Value = '  +2.47%  '
$ws.Range("E12").Value = '  +0.57%  '
$ws.Range("D13").Value = '2.860.36'
$ws.Range("E14").Value = '  +1.39%  '
$ws.Range("D15").Value = '15.72'
$ws.Range("E15").Value = '  -0.50%  '
$ws.Range("D16").Value = '2.467.68'
$ws.Range("E16").Value = '  +0.41%  '
$ws.Range("E17").Value = '  +3.91%  '
$ws.Range("D18").Value = '41.751.57'
$ws.Range("E18").Value = '  +0.53%  '
$ws.Range("E19").Value = '  +0.12%  '
$ws.Range("E20").Value = '  +1.71%  '
$ws.Range("E21").Value = '  +0.50%  '
$ws.Range("D22").Value = '11.38'
$ws.Range("E22").Value = '  +2.23%  '
$ws.Range("E23").Value = '  +1.62%  '
$ws.Range("E24").Value = '  +1.60%  '
$ws.Range("E25").Value = '  +3.16%  '
$ws.Range("E26").Value = '  +0.01%  '
$ws.Range("D27").Value = '25.26'
$ws.Range("E27").Value = '  +2.88%  '
$ws.Range("E28").Value = '  +0.77%  '
$ws.Range("E29").Value = '  +1.12%  '
$ws.Range("D30").Value = '37.15'
$ws.Range("E30").Value = '  +5.33%  '
$ws.Range("D31").Value = '159.42'
$ws.Range("E31").Value = '  +2.07%  '
$ws.Range("E32").Value = '  +2.03%  '
$ws.Range("E33").Value = '  -0.01%  '
$ws.Range("E34").Value = '  +1.24%  '
$ws.Range("E35").Value = '  -0.42%  '
$ws.Range("D36").Value = '17.44'
$ws.Range("E36").Value = '  +1.28%  '
$ws.Range("E37").Value = '  +5.28%  '
$ws.Range("E39").Value = '  +1.82%  '
$ws.Range("D40").Value = '0.104'
$ws.Range("E40").Value = '  +1.45%  '
$ws.Range("E41").Value = '  +1.38%  '
$ws.Range("D42").Value = '2.52'
$ws.Range("E42").Value = '  +7.73%  '
$ws.Range("D43").Value = '2.011.10'
$ws.Range("E43").Value = '  +3.56%  '
$ws.Range("D44").Value = '19.28'
$ws.Range("E44").Value = '  +2.57%  '
$ws.Range("E45").Value = '  +0.99%  '
$ws.Range("E46").Value = '  +3.30%  '
$ws.Range("D47").Value = '9.49'
$ws.Range("E47").Value = '  +5.14%  '
$ws.Range("D48").Value = '2.716.33'
$ws.Range("E48").Value = '  +0.45%  '
$ws.Range("E49").Value = '  +7.97%  '
$ws.Range("D50").Value = '98.36'
$ws.Range("E50").Value = '  +1.78%  '
$ws.Range("D51").Value = '67.45'
$ws.Range("E51").Value = '  +0.80%  '
